$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header cells F1/G1 (C1 -> C5, C2 -> C6)
$ws.Range("F1").Value = "C5"
$ws.Range("G1").Value = "C6"

# Add new header cells I1/J1, matching the style of the other header cells (copy from H1)
$ws.Range("H1").Copy($ws.Range("I1:J1"))

$ws.Range("I1").Value = "C1"
$ws.Range("J1").Value = "C2"

# Fill in new data columns I and J for rows 2-7
$ws.Range("I2").Value = "C4"
$ws.Range("J2").Value = "C3"

$ws.Range("I3").Value = "C4"
$ws.Range("J3").Value = "C3"

$ws.Range("I4").Value = "C4"
$ws.Range("J4").Value = "C3"

$ws.Range("I5").Value = "C4"
$ws.Range("J5").Value = "C3"

$ws.Range("I6").Value = "C1"
$ws.Range("J6").Value = "C2"

$ws.Range("I7").Value = "C1"
$ws.Range("J7").Value = "C2"
